# Apply the Alvearie FHIR IG regeneration update:
#  - Metadata sheet: version bump, new date, new publisher, replace the
#    duplicated "Contact" rows with a single "Jurisdiction" row, which
#    removes one row overall (21 -> 20 rows).
#  - Elements sheet: the Extension row's "Short"/"Definition" columns now
#    carry the StructureDefinition's own title/description text.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Remove the duplicated second "Contact" / "No display for ContactDetail"
# row (originally row 11); this shifts rows 12-21 up to 11-20, matching
# the new A1:B20 used range.
$meta.Rows.Item(11).Delete()

# Version bump
$meta.Cells.Item(3, 2).Value = "6.0.0"

# Publication date refresh
$meta.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher now populated
$meta.Cells.Item(9, 2).Value = "Alvearie Team"

# Former "Contact" row becomes "Jurisdiction"
$meta.Cells.Item(10, 1).Value = "Jurisdiction"
$meta.Cells.Item(10, 2).Value = "United States of America"

$elements = $wb.Worksheets.Item("Elements")

# Row 2 (the root "Extension" element) now shows the profile's own
# Short/Definition text instead of the generic Extension boilerplate.
$elements.Cells.Item(2, 11).Value = "Longterm Care Waiting Period"
$elements.Cells.Item(2, 12).Value = "Elimination or waiting period for the long-term disability (LTD) benefit (for example, 90, 180 or 365 days). This is the amount of time between the first absent date and the coverage begin date."
